$wb = $excel.ActiveWorkbook

# New file name (markdown source doc) that replaces the old UUID-based name everywhere.
$newFileName = "e33ee5b4-30a9-4314-baf0-637d772dabb8.md"
$newFileUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/73bc401113958c29f19d27afca921ba4fd1053f4/e2e/$newFileName"
$configUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/73bc401113958c29f19d27afca921ba4fd1053f4/.localization-config"

$newStatus = "Handoff failed"
$epoch     = "0001-01-01 00:00:00"
$ignored   = "Ignored"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.UsedRange.Hyperlinks.Delete()
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $newFileUrl, "", "", $newFileName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config")

# ---- zh-cn sheet ----
$wsZh.UsedRange.Hyperlinks.Delete()
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Clear()
$wsZh.Range("D2").Value = $epoch
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = $ignored
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $newFileUrl, "", "", $newFileName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, "", "", ".localization-config")

# ---- de-de sheet ----
$wsDe.UsedRange.Hyperlinks.Delete()
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Clear()
$wsDe.Range("D2").Value = $epoch
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = $ignored
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $newFileUrl, "", "", $newFileName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, "", "", ".localization-config")
